# Generate Report for Archive
# The handoff status for the files "6efe1dd5-...md" (row 8) and
# "df976f6e-...md" (row 9) moved from "Ready for handoff" to "In Translation"
# on the Overview sheet (zh-cn / de-de status columns) as well as on each of
# the per-locale sheets (Status column). The third file in this batch,
# "e9ebefe9-...md" (row 10), stays "Ready for handoff".

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B8").Value = "In Translation"
$ws.Range("C8").Value = "In Translation"
$ws.Range("B9").Value = "In Translation"
$ws.Range("C9").Value = "In Translation"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C8").Value = "In Translation"
$wsZh.Range("C9").Value = "In Translation"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C8").Value = "In Translation"
$wsDe.Range("C9").Value = "In Translation"
